$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the travel name for the "bridge" entry to match the new wording
$ws.Range("B9").Value = "Cross a bridge"

# Move the active selection to C14, matching the edited sheet view state
$ws.Range("C14").Select()
